$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optional indicators")
$ws.Range("A100").Value = "test"
